$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B's new value ("102") is a purely numeric-looking string that must
# be stored as text (it becomes its own shared string in the target file,
# not a numeric cell). Mark the cells as Text before writing so the engine
# keeps the value as a string instead of auto-detecting it as a number,
# then drop the number formatting again so the cells end up using the
# sheet's normal/default style, same as every other cell around them.
$ws.Range("B5:B6").NumberFormat = "@"

$ws.Range("A5").Value = "od-856429"
$ws.Range("B5").Value = "102"
$ws.Range("C5").Value = "user101"
$ws.Range("D5").Value = "2019-Sep-13"

$ws.Range("A6").Value = "od-856429"
$ws.Range("B6").Value = "102"
$ws.Range("C6").Value = "user101"
$ws.Range("D6").Value = "2019-Sep-13"

$ws.Range("B5:B6").ClearFormats()
